$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values of rows 2, 3, 4 for the columns that change.
$cols = @("A","B","E","F","G","H","Q","R")

$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value2
    $row3[$col] = $ws.Range("${col}3").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
}

# Target mapping (cyclic rotation): new row2 = old row4, new row3 = old row2, new row4 = old row3
foreach ($col in $cols) {
    $ws.Range("${col}2").Value2 = $row4[$col]
    $ws.Range("${col}3").Value2 = $row2[$col]
    $ws.Range("${col}4").Value2 = $row3[$col]
}
